$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.197.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.63%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.112.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.60%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.19%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.72'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.36%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.104.22'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.55%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.38%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +13.90%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +7.90%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.68%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +8.55%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.61%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.64%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.626.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.56%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.17'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.86%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.114.61'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.41%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.107.68'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.19'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.54%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.20'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.01%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.28%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.32'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.81'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.56%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.48'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +9.86%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.21%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.41%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.03%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.86'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +10.53%  '

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.111'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.16%  '

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.93'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.97%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0872'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +14.40%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +16.41%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.37%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.07'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.17%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.31'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +19.93%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.78'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.77%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '438.68'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.74'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.75%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.915.96'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.24%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.07%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +12.70%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +7.23%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +8.42%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.11'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.70%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.00%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.65'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.10%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.55'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.54%  '
